# APPMO-SP_CRE version bump: 2.0 -> 1.2 in the "CONTROL DE VERSIONES"
# table, plus relocating the hidden "_GoBack" bookmark from its old
# spot (the meeting-date cell, between "11" and "/Agosto") to the new
# last-edit spot (right after the new "1.2" text) - exactly what Word
# itself does automatically to track the most recent edit position.
#
# NOTE: table/cell/range references are always re-fetched from
# $d right before use (never reused across an edit) because this
# headless engine does not keep older references in sync with
# subsequent document mutations.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Update the version number in the "CONTROL DE VERSIONES" table
# ---------------------------------------------------------------
$d.Tables.Item(1).Cell(3, 1).Range.Find.Execute("2.0", $true, $false, $false, `
    $false, $false, $true, 1, $false, "1.2", 2) | Out-Null

# ---------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark, which currently sits between
#    the "11" and "/Agosto" runs of the last table's date cell
#    ("11/Agosto/2019"). This engine only tears a bookmark down when
#    an edit actually straddles its (zero-width) position, and only
#    treats an edit as "real" when the replacement text differs from
#    what was there before - so nudge the text, confirm the bookmark
#    is gone, then put the original characters back.
# ---------------------------------------------------------------
$dateSearch = $d.Tables.Item(3).Cell(13, 4).Range.Duplicate
$dateSearch.Find.Execute("11/Agosto") | Out-Null
$anchorStart = $dateSearch.Start

# The old bookmark boundary sits right between "11" and "/Agosto".
$boundary = $anchorStart + 2

$crossRng = $d.Range($boundary - 1, $boundary + 1)
$originalText = $crossRng.Text
$crossRng.Text = "##"
$crossRng2 = $d.Range($boundary - 1, $boundary + 1)
$crossRng2.Text = $originalText

# The text replace above merges the three runs that made up this
# cell's paragraph ("11" / "/Agosto" / "/2019") into a single run;
# split it back apart at the original boundaries by toggling a
# character property on/off over each sub-range (this engine -
# like Word - always breaks a run apart when formatting is applied
# to only part of it).
$afterSlashAgosto = $boundary + 7
$runSplit1 = $d.Range($boundary, $afterSlashAgosto)
$runSplit1.Bold = 1
$runSplit1.Bold = 0
$runSplit2 = $d.Range($afterSlashAgosto, $afterSlashAgosto + 5)
$runSplit2.Bold = 1
$runSplit2.Bold = 0

# ---------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark right after the new "1.2" text.
#    A collapsed (zero-length) Range cannot be handed straight to
#    Bookmarks.Add in this engine (it silently resets to 0,0), so
#    bracket a throwaway character right after "1.2", bookmark that
#    1-character range, then delete the character - the bookmark
#    naturally collapses to zero width in the correct spot, just
#    like Word leaves it. (Collapse() on a Range also does not feed
#    through correctly into a later InsertBefore/InsertAfter call in
#    this engine, so insert straight off of the still-uncollapsed
#    "1.2" hit range instead.)
# ---------------------------------------------------------------
$versionFind = $d.Tables.Item(1).Cell(3, 1).Range.Duplicate
$versionFind.Find.Execute("1.2") | Out-Null
$afterVersion = $versionFind.End
$versionFind.InsertAfter("#")

$tempCharRng = $d.Range($afterVersion, $afterVersion + 1)
$d.Bookmarks.Add("_GoBack", $tempCharRng) | Out-Null

$tempCharRng2 = $d.Range($afterVersion, $afterVersion + 1)
$tempCharRng2.Delete()
